$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '49.153.24'
$ws.Range("E2").Value = '  -1.27%  '
$ws.Range("D3").Value = '2.630.08'
$ws.Range("E3").Value = '  +0.38%  '
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = '  +0.01%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '111.41'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +1.02%  '
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '322.66'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  -0.36%  '
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.525'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  -1.89%  '
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -3.34%  '
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.76'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  -2.84%  '
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.82'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  -3.95%  '
$ws.Range("E12").Value = '  -1.72%  '
$ws.Range("E13").Value = '  +0.24%  '
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.25'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  -1.09%  '
$ws.Range("D15").Value = '3.044.74'
$ws.Range("E15").Value = '  +0.46%  '
$ws.Range("D16").Value = '2.633.28'
$ws.Range("E16").Value = '  +1.03%  '
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.860'
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  -1.41%  '
$ws.Range("D18").Value = '49.135.16'
$ws.Range("E18").Value = '  -1.34%  '
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.98'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  -4.14%  '
$ws.Range("E20").Value = '  -3.55%  '
$ws.Range("E21").Value = '  -1.26%  '
$ws.Range("D22").Value = '0.0₃0945'
$ws.Range("E22").Value = '  -0.98%  '
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '269.40'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  -4.28%  '
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.50'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  -6.00%  '
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.53'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  -1.39%  '
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.10'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  -2.22%  '
$ws.Range("E27").Value = '  -0.07%  '
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.03'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  +0.57%  '
$ws.Range("E29").Value = '  -0.89%  '
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.09'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  -3.45%  '
$ws.Range("E31").Value = '  -4.70%  '
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.52'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  -0.18%  '
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.48'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  +0.15%  '
$ws.Range("E34").Value = '  -0.31%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0798'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  +0.19%  '
$ws.Range("B36").Value = 'Celestia'
$ws.Range("C36").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.02'
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = '  -3.41%  '
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.98'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  +4.35%  '
$ws.Range("E38").Value = '  -0.93%  '
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.13'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  +1.23%  '
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '127.10'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  +2.96%  '
$ws.Range("E41").Value = '  -1.77%  '
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.06'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  -3.82%  '
$ws.Range("E43").Value = '  -4.40%  '
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0316'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").Value = '2.065.17'
$ws.Range("E45").Value = '  +0.40%  '
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.15'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  +5.90%  '
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.25'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  -3.20%  '
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.13'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  -2.85%  '
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.87'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  -2.14%  '
$ws.Range("E50").Value = '  -3.62%  '
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '58.58'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  +1.32%  '
